$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format in bulk first, to avoid Excel
# auto-converting numeric-looking strings (e.g. "1.00", "0.999") into numbers
# and stripping their literal formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.560.25'
$ws.Range("E2").Value = '  -0.98%  '

$ws.Range("D3").Value = '3.073.39'
$ws.Range("E3").Value = '  -0.80%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.46%  '

$ws.Range("D5").Value = '591.59'
$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("D6").Value = '154.43'
$ws.Range("E6").Value = '  +0.49%  '

$ws.Range("E7").Value = '  -0.27%  '

$ws.Range("D8").Value = '0.539'
$ws.Range("E8").Value = '  +1.22%  '

$ws.Range("D9").Value = '3.070.46'
$ws.Range("E9").Value = '  -0.79%  '

$ws.Range("E10").Value = '  -1.19%  '

$ws.Range("D11").Value = '5.91'
$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("D12").Value = '0.452'
$ws.Range("E12").Value = '  -1.58%  '

$ws.Range("E13").Value = '  -2.43%  '

$ws.Range("D14").Value = '36.83'
$ws.Range("E14").Value = '  -2.17%  '

$ws.Range("E15").Value = '  +1.27%  '

$ws.Range("D16").Value = '3.575.15'
$ws.Range("E16").Value = '  -1.09%  '

$ws.Range("D17").Value = '7.20'
$ws.Range("E17").Value = '  +0.30%  '

$ws.Range("D18").Value = '63.463.25'
$ws.Range("E18").Value = '  -0.73%  '

$ws.Range("D19").Value = '3.065.51'
$ws.Range("E19").Value = '  -0.99%  '

$ws.Range("D20").Value = '484.25'
$ws.Range("E20").Value = '  +2.89%  '

$ws.Range("D21").Value = '14.61'
$ws.Range("E21").Value = '  -0.39%  '

$ws.Range("D22").Value = '0.710'
$ws.Range("E22").Value = '  -3.11%  '

$ws.Range("D23").Value = '7.59'
$ws.Range("E23").Value = '  +0.61%  '

$ws.Range("D24").Value = '2.42'
$ws.Range("E24").Value = '  +1.75%  '

$ws.Range("D25").Value = '82.01'
$ws.Range("E25").Value = '  +0.66%  '

$ws.Range("D26").Value = '12.98'
$ws.Range("E26").Value = '  -1.91%  '

$ws.Range("E27").Value = '  +6.72%  '

$ws.Range("E28").Value = '  +0.26%  '

$ws.Range("D29").Value = '7.52'
$ws.Range("E29").Value = '  +2.02%  '

$ws.Range("D30").Value = '2.24'
$ws.Range("E30").Value = '  +2.15%  '

$ws.Range("D31").Value = '2.69'
$ws.Range("E31").Value = '  -0.33%  '

$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.57%  '

$ws.Range("D33").Value = '27.37'
$ws.Range("E33").Value = '  -0.15%  '

$ws.Range("E34").Value = '  -2.95%  '

$ws.Range("D35").Value = '1.07'
$ws.Range("E35").Value = '  +1.29%  '

$ws.Range("E36").Value = '  -3.53%  '

$ws.Range("D37").Value = '6.07'
$ws.Range("E37").Value = '  -1.25%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = '2.23'
$ws.Range("E38").Value = '  -1.18%  '

$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").Value = '3.25'
$ws.Range("E39").Value = '  -4.43%  '

$ws.Range("D40").Value = '9.34'
$ws.Range("E40").Value = '  +0.56%  '

$ws.Range("D41").Value = '50.65'
$ws.Range("E41").Value = '  +0.13%  '

$ws.Range("D42").Value = '440.67'
$ws.Range("E42").Value = '  -2.25%  '

$ws.Range("D43").Value = '0.289'
$ws.Range("E43").Value = '  +0.43%  '

$ws.Range("E44").Value = '  +3.01%  '

$ws.Range("D45").Value = '0.0365'
$ws.Range("E45").Value = '  -0.84%  '

$ws.Range("D46").Value = '2.831.50'
$ws.Range("E46").Value = '  -0.27%  '

$ws.Range("D47").Value = '39.69'
$ws.Range("E47").Value = '  +0.24%  '

$ws.Range("D48").Value = '132.87'
$ws.Range("E48").Value = '  +2.14%  '

$ws.Range("D49").Value = '25.59'
$ws.Range("E49").Value = '  +0.45%  '

$ws.Range("D51").Value = '2.24'
$ws.Range("E51").Value = '  -0.63%  '

# Restore default style on column D so no stray number-format styling remains
$ws.Range("D2:D51").Style = "Normal"
